$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 5469
$ws.Range("A2").Value = 4783.999999999956
$ws.Range("A3").Value = 4130
$ws.Range("A4").Value = 4669
$ws.Range("A5").Value = 4914.999999999997
$ws.Range("A6").Value = 5160
$ws.Range("A7").Value = 5100.999999999956
$ws.Range("A8").Value = 4996.999999999985
